# Script applies three changes to Sheet1:
#  1. Swap the match data (columns F:V) between rows 4 and 5.
#  2. Swap the match data (columns F:V) between rows 62 and 63.
#  3. Append a new match row (row 64): Esteghlal Khuzestan vs Gol Gohar,
#     copying the formatting of the row above (row 63) for the styled
#     columns A and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap rows 4 and 5 (columns F:V only; A:E - index/meta/date - stay put) ---
$row4 = $ws.Range("F4:V4").Value2
$row5 = $ws.Range("F5:V5").Value2
$ws.Range("F4:V4").Value2 = $row5
$ws.Range("F5:V5").Value2 = $row4

# --- 2. Swap rows 62 and 63 (columns F:V only) ---
$row62 = $ws.Range("F62:V62").Value2
$row63 = $ws.Range("F63:V63").Value2
$ws.Range("F62:V62").Value2 = $row63
$ws.Range("F63:V63").Value2 = $row62

# --- 3. Add new row 64 ---
# Copy formatting from row 63 for the styled cells (A: bold/border/centered
# index style, E: date-time number format) so the new row matches existing
# sheet conventions, then set the values/content for every column.
$ws.Range("A63").Copy() | Out-Null
$ws.Range("A64").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$ws.Range("E63").Copy() | Out-Null
$ws.Range("E64").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A64").Value2 = 63
$ws.Range("B64").Value2 = "iran"
$ws.Range("C64").Value2 = "persian-gulf-pro-league"
$ws.Range("D64").Value2 = "2023-2024"
$ws.Range("E64").Value2 = 45233.64583333334
$ws.Range("F64").Value2 = "Esteghlal Khuzestan"
$ws.Range("G64").Value2 = 1
$ws.Range("H64").Value2 = "Gol Gohar"
$ws.Range("I64").Value2 = 1
$ws.Range("J64").Value2 = 2.76
$ws.Range("K64").Value2 = "02/11/2023 03:42"
$ws.Range("L64").Value2 = 5.79
$ws.Range("M64").Value2 = "03/11/2023 15:29"
$ws.Range("N64").Value2 = 2.65
$ws.Range("O64").Value2 = "02/11/2023 03:42"
$ws.Range("P64").Value2 = 2.58
$ws.Range("Q64").Value2 = "03/11/2023 15:29"
$ws.Range("R64").Value2 = 2.68
$ws.Range("S64").Value2 = "02/11/2023 03:42"
$ws.Range("T64").Value2 = 1.95
$ws.Range("U64").Value2 = "03/11/2023 15:29"
$ws.Range("V64").Value2 = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-khuzestan-gol-gohar/M3EzdSgG/"
